$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 45
$ws.Range("F3").Value = 26688
$ws.Range("F4").Value = 589
$ws.Range("F6").Value = 610
$ws.Range("F7").Value = 177
$ws.Range("F8").Value = 552
$ws.Range("F10").Value = 360
$ws.Range("F11").Value = 242
$ws.Range("F12").Value = 190
$ws.Range("F15").Value = 65
$ws.Range("F16").Value = 426
$ws.Range("F17").Value = 59
$ws.Range("F18").Value = 1544
$ws.Range("F19").Value = 209
$ws.Range("F20").Value = 45
$ws.Range("F21").Value = 438

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 232
$ws.Range("F6").Value = 198
$ws.Range("F7").Value = 198
$ws.Range("F11").Value = 440
$ws.Range("F16").Value = 60

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5067
$ws.Range("F3").Value = 228

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 45
$ws.Range("F3").Value = 5067
$ws.Range("F4").Value = 228
$ws.Range("F5").Value = 26688
$ws.Range("F6").Value = 589
$ws.Range("F9").Value = 232
$ws.Range("F10").Value = 610
$ws.Range("F13").Value = 177
$ws.Range("F14").Value = 198
$ws.Range("F15").Value = 198
$ws.Range("F19").Value = 440
$ws.Range("F20").Value = 552
$ws.Range("F23").Value = 360
$ws.Range("F24").Value = 242
$ws.Range("F25").Value = 190
$ws.Range("F29").Value = 65
$ws.Range("F32").Value = 426
$ws.Range("F33").Value = 59
$ws.Range("F34").Value = 60
$ws.Range("F35").Value = 1544
$ws.Range("F36").Value = 209
$ws.Range("F38").Value = 45
$ws.Range("F39").Value = 438
